$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B5").Value = "حلىمى علي ريان"
$ws.Range("B6").Value = "حسين محمد ماهر براء الدين"
$ws.Range("B7").Value = "محمد منير عبد الحميد كمال"
$ws.Range("B22").Value = "عمرو أيهاب مختار قرخات"
$ws.Range("B32").Value = "محمد عبد الرحيم سالم عبد الهادي"
$ws.Range("B34").Value = "زياد ايهاه محمد ممدوح ناقع"
$ws.Range("B35").Value = "اية احمد محمد خشبه"
$ws.Range("B38").Value = "جنى ايمن وقائى محمد عيسى"
$ws.Range("B39").Value = "سلمى محمد ابراهيم قتحدى ابوريدة"
$ws.Range("B40").Value = "مروان عمرو عبد المجيد فؤاد احمد شكرى"
$ws.Range("B47").Value = "جمانة عمرو مصطفى عبد الصالح عرابي"
$ws.Range("B51").Value = "حمزة احمد محمد منير الجوهري"
